# Budweiser Case Study deck cleanup:
#  1. Fix a typo ("AVB" -> "ABV") on the "Median ABV vs Median IBU" slide.
#  2. Remove three now-redundant/obsolete slides whose content duplicates
#     (a less complete version of) material kept elsewhere in the deck:
#       - "States with highest ABV"            (slide id 265)
#       - "Summary Statistics of ABV"           (slide id 262)
#       - "Relationship between IBU and ABV"    (slide id 267)

$p = $ppt.ActivePresentation

function Get-SlideTitleText($slide) {
    foreach ($sh in $slide.Shapes) {
        if ($sh.Type -eq 14 -and $sh.HasTextFrame) {
            return $sh.TextFrame.TextRange.Text
        }
    }
    return $null
}

# --- 1. Typo fix: "AVB" -> "ABV" on the "Median ABV vs Median IBU" slide ---
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    if ((Get-SlideTitleText $slide) -eq "Median ABV vs Median IBU") {
        $contentShape = $slide.Shapes.Item("Content Placeholder 2")
        $para1 = $contentShape.TextFrame.TextRange.Paragraphs(1)
        $pos = $para1.Text.IndexOf("AVB")
        if ($pos -ge 0) {
            $badWord = $para1.Characters($pos + 1, 3)
            $badWord.Text = "ABV"
        }
        break
    }
}

# --- 2. Delete the obsolete/duplicate slides --------------------------------
$titlesToRemove = @(
    "States with highest ABV",
    "Summary Statistics of ABV",
    "Relationship between IBU and ABV"
)

foreach ($title in $titlesToRemove) {
    for ($i = 1; $i -le $p.Slides.Count; $i++) {
        $slide = $p.Slides.Item($i)
        if ((Get-SlideTitleText $slide) -eq $title) {
            $slide.Delete()
            break
        }
    }
}
